$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original A18 formatting (date-style) before it gets overwritten,
# so the brand-new row 19 can reuse the same style index instead of a freshly
# minted one.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Rewrite the data rows bottom-up so each row ends up holding the values the
# diff shows for that row (effectively an insert-at-row-2 + shift, done via
# plain value writes so no incidental style churn is introduced).
$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = -0.6470065423293869
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 1.000784331770643

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = -0.2814561130375925
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = 0.06837453075889677

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 3.293290997728171
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = 1.50842747477733

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = 0.4255262881966981
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 2.107044281672787

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = -3.258619210312885
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = -2.035749133083187

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 1.457852003181337
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 0.7756897792100315

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 1.911050033324102
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = 2.362613045431528

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 2.18621550610123
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 2.357704431248409

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 2.204449574611278
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = 1.401113624217043

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 1.365576377841027
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 1.582979977679533

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 1.297015177357297
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 0.309136676902555

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = 0.2545814083968478
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 1.223618887196531

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 0.4451370000809973
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 0.6480763427741953

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 2.791140000794279
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 1.722110645261998

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = 0.9337833426867448
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 1.899713704008654

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = -0.8792832172735965
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = -0.10374497415091

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 1.381024225294869
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 1.929644353290927

$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 1.144978573787081
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 1.854677717601594
